# Update Sage scrape results
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rsquo = [char]0x2019

$ws.Range("B2").Value = "The Impact of Cyber Conflict on International Interactions"

$ws.Range("B3").Value = "Evaluation of communIT, a large-scale, cyber-physical artifact supporting diverse subgroups building community"
$ws.Range("F3").Value = "Restricted"

$ws.Range("B4").Value = "Securing Virtual Space: Cyber War, Cyber Terror, and Risk"
$ws.Range("F4").Value = "Restricted"

$ws.Range("B5").Value = "Warring from the virtual to the real: Assessing the public" + $rsquo + "s threshold for war over cyber security"

$ws.Range("B6").Value = "Accountability and cyber conflict: examining institutional constraints on the use of cyber proxies"

$ws.Range("B7").Value = "Adaptive Torque and Position Control for a Legged Robot Based on a Series Elastic Actuator"
$ws.Range("F7").Value = "Restricted"

$ws.Range("B8").Value = "Modeling and application for pneumatic soft actuators based on a novel deep neural network"
$ws.Range("F8").Value = "Restricted"

$ws.Range("B9").Value = "On 3D simultaneous attack against manoeuvring target with communication delays"
$ws.Range("F9").Value = "Restricted"

$ws.Range("B10").Value = "Ontological security, cyber technology, and states" + $rsquo + " responses"
$ws.Range("F10").Value = "Restricted"

$ws.Range("B11").Value = "Responding to Uncertainty: The Importance of Covertness in Support for Retaliation to Cyber and Kinetic Attacks"
